# Add BEFORE/AFTER keyword support to the rule expression in C2, and move
# the sheet's active selection from B2 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C2")
$text = $cell.Value()

# The rule text contains three "} AND {" joins. The first one (between the
# token/regex condition and the first "Node down" condition) gets an
# "AFTER 00:00:00:980" qualifier; the second one (between the first and
# second "Node down" conditions) gets a "BEFORE 00:00:00:876" qualifier.
$idx1 = $text.IndexOf('} AND {')
$text = $text.Substring(0, $idx1) + '} AND AFTER 00:00:00:980 {' + $text.Substring($idx1 + '} AND {'.Length)

$idx2 = $text.IndexOf('} AND {', $idx1 + 1)
$text = $text.Substring(0, $idx2) + '} AND BEFORE 00:00:00:876 {' + $text.Substring($idx2 + '} AND {'.Length)

$cell.Value = $text

# Move the active selection from B2 to A2.
$ws.Range("A2").Select()
